# Add a "Typing" row to the "Flight Mission Cycle" summary sheet, mirroring
# the existing rows that summarize the number of cycles for each mode.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flight Mission Cycle")

$ws.Range("A4").Value = "Typing"
$ws.Range("B4").Value = 2

# Make this sheet the active one, with I7 selected, matching the edited file.
$ws.Activate()
$ws.Range("I7").Select()
